$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'65.935.74"
$ws.Range("E2").Value2 = "  +3.58%  "
$ws.Range("D3").Value2 = "'3.417.84"
$ws.Range("E3").Value2 = "  +2.41%  "
$ws.Range("D4").Value2 = "'0.999"
$ws.Range("E4").Value2 = "  +0.09%  "
$ws.Range("D5").Value2 = "'563.00"
$ws.Range("E5").Value2 = "  +2.57%  "
$ws.Range("D6").Value2 = "'175.84"
$ws.Range("E6").Value2 = "  +2.45%  "
$ws.Range("E7").Value2 = "  +2.95%  "
$ws.Range("D8").Value2 = "'3.411.78"
$ws.Range("E8").Value2 = "  +2.51%  "
$ws.Range("E9").Value2 = "  +0.00%  "
$ws.Range("E10").Value2 = "  +14.55%  "
$ws.Range("D11").Value2 = "'0.634"
$ws.Range("E11").Value2 = "  +3.37%  "
$ws.Range("D12").Value2 = "'55.18"
$ws.Range("E12").Value2 = "  +3.49%  "
$ws.Range("D13").Value2 = "'0.0000283"
$ws.Range("E13").Value2 = "  +6.68%  "
$ws.Range("D14").Value2 = "'9.18"
$ws.Range("E14").Value2 = "  +3.06%  "
$ws.Range("D15").Value2 = "'3.948.39"
$ws.Range("E15").Value2 = "  +2.43%  "
$ws.Range("D16").Value2 = "'18.37"
$ws.Range("E16").Value2 = "  +3.04%  "
$ws.Range("D17").Value2 = "'3.398.19"
$ws.Range("E17").Value2 = "  +2.69%  "
$ws.Range("D18").Value2 = "'0.119"
$ws.Range("E18").Value2 = "  +1.68%  "
$ws.Range("D19").Value2 = "'65.675.93"
$ws.Range("E19").Value2 = "  +3.50%  "
$ws.Range("D20").Value2 = "'11.92"
$ws.Range("E20").Value2 = "  +2.09%  "
$ws.Range("D21").Value2 = "'0.995"
$ws.Range("E21").Value2 = "  +2.52%  "
$ws.Range("D22").Value2 = "'471.79"
$ws.Range("E22").Value2 = "  +15.17%  "
$ws.Range("D23").Value2 = "'5.18"
$ws.Range("E23").Value2 = "  +19.63%  "
$ws.Range("D24").Value2 = "'4.15"
$ws.Range("E24").Value2 = "  +2.37%  "
$ws.Range("D25").Value2 = "'86.69"
$ws.Range("E25").Value2 = "  +4.49%  "
$ws.Range("D26").Value2 = "'13.54"
$ws.Range("E26").Value2 = "  +2.22%  "
$ws.Range("D27").Value2 = "'10.94"
$ws.Range("E27").Value2 = "  +3.50%  "
$ws.Range("D28").Value2 = "'2.90"
$ws.Range("E28").Value2 = "  +6.66%  "
$ws.Range("D29").Value2 = "'8.93"
$ws.Range("E29").Value2 = "  +4.23%  "
$ws.Range("D30").Value2 = "'31.07"
$ws.Range("E30").Value2 = "  +6.75%  "
$ws.Range("D31").Value2 = "'6.72"
$ws.Range("E31").Value2 = "  +5.23%  "
$ws.Range("B32").Value2 = "OKB"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value2 = "'63.38"
$ws.Range("E32").Value2 = "  +10.43%  "
$ws.Range("B33").Value2 = "Cosmos"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value2 = "'11.58"
$ws.Range("E33").Value2 = "  +2.41%  "
$ws.Range("D34").Value2 = "'583.30"
$ws.Range("E34").Value2 = "  +2.09%  "
$ws.Range("E35").Value2 = "  +2.69%  "
$ws.Range("E36").Value2 = "  +0.19%  "
$ws.Range("E37").Value2 = "  -3.13%  "
$ws.Range("E38").Value2 = "  +3.34%  "
$ws.Range("B39").Value2 = "PEPE"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value2 = "'0.0₃0762"
$ws.Range("E39").Value2 = "  +3.15%  "
$ws.Range("B40").Value2 = "InjectiveProtocol"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value2 = "'35.92"
$ws.Range("E40").Value2 = "  +2.00%  "
$ws.Range("D41").Value2 = "'0.375"
$ws.Range("E41").Value2 = "  +2.34%  "
$ws.Range("D42").Value2 = "'3.099.35"
$ws.Range("E42").Value2 = "  -1.64%  "
$ws.Range("D43").Value2 = "'0.997"
$ws.Range("E43").Value2 = "  +0.08%  "
$ws.Range("D44").Value2 = "'2.85"
$ws.Range("E44").Value2 = "  +1.18%  "
$ws.Range("D45").Value2 = "'0.0417"
$ws.Range("E45").Value2 = "  +3.87%  "
$ws.Range("B46").Value2 = "Fetch.AI"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value2 = "'2.50"
$ws.Range("E46").Value2 = "  +3.44%  "
$ws.Range("B47").Value2 = "ApeXProtocol"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value2 = "'3.22"
$ws.Range("E47").Value2 = "  -0.90%  "
$ws.Range("E48").Value2 = "  +5.87%  "
$ws.Range("E49").Value2 = "  -0.69%  "
$ws.Range("D50").Value2 = "'8.40"
$ws.Range("E50").Value2 = "  +5.31%  "
$ws.Range("D51").Value2 = "'136.60"
$ws.Range("E51").Value2 = "  +2.90%  "
